$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the Status column (E2:E17) values to "AVAILABLE"
$ws.Range("E2:E17").Value = "AVAILABLE"

# Update the selection to reflect the new range E2:E17
$ws.Range("E2:E17").Select()
